$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (source, amount, date-serial)
$data = @(
    @("misc",        126000, 46011.22928240741),
    @("shop",         30000, 46009.22928240741),
    @("dedef",         2000, 46009.22928240741),
    @("New",          20000, 46001.22928240741),
    @("hello",        30000, 46001.22928240741),
    @("hlkjfdnvjkn",  20000, 45996.22928240741),
    @("Salary",      200000, 45992.22928240741),
    @("project 2",    10000, 45992.22928240741)
)

# Use the existing date-formatted cell (C2) as the style template for every
# date cell so newly created rows pick up the same number format (style s="1").
$ws.Range("C2").Copy()

$row = 2
foreach ($entry in $data) {
    $ws.Range("A" + $row).Value = $entry[0]
    $ws.Range("B" + $row).Value = $entry[1]
    $ws.Range("C" + $row).PasteSpecial(-4122)
    $ws.Range("C" + $row).Value = $entry[2]
    $row = $row + 1
}
